$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '57.829.10'
$ws.Range('E2').Value = '  -0.37%  '

$ws.Range('D3').Value = '3.105.14'
$ws.Range('E3').Value = '  +1.22%  '

$ws.Range('E4').Value = '  +0.00%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '524.93'

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '142.19'
$ws.Range('E6').Value = '  +0.22%  '

$ws.Range('E7').Value = '  -0.05%  '

$ws.Range('D8').Value = '3.105.57'
$ws.Range('E8').Value = '  +1.30%  '

$ws.Range('E9').Value = '  +0.71%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '7.23'
$ws.Range('E10').Value = '  -0.92%  '

$ws.Range('E11').Value = '  +0.38%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.389'
$ws.Range('E12').Value = '  +3.53%  '

$ws.Range('D13').Value = '3.636.51'
$ws.Range('E13').Value = '  +1.22%  '

$ws.Range('E14').Value = '  +1.33%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '25.81'
$ws.Range('E15').Value = '  -2.26%  '

$ws.Range('E16').Value = '  +0.19%  '

$ws.Range('D17').Value = '57.896.54'
$ws.Range('E17').Value = '  -0.27%  '

$ws.Range('D18').Value = '3.104.36'
$ws.Range('E18').Value = '  +1.27%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.11'
$ws.Range('E19').Value = '  -0.13%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '12.81'
$ws.Range('E20').Value = '  -0.18%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '8.04'
$ws.Range('E21').Value = '  -1.68%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '341.70'
$ws.Range('E22').Value = '  +2.93%  '

$ws.Range('E23').Value = '  +0.16%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.513'
$ws.Range('E24').Value = '  +2.27%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '66.97'
$ws.Range('E25').Value = '  +2.37%  '

$ws.Range('E26').Value = '  -0.41%  '

$ws.Range('E27').Value = '  -0.11%  '

$ws.Range('D28').Value = '0.0₃0922'
$ws.Range('E28').Value = '  +1.25%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '6.51'
$ws.Range('E29').Value = '  +0.35%  '

$ws.Range('E30').Value = '  +0.08%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '7.21'
$ws.Range('E31').Value = '  -0.51%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.88'
$ws.Range('E32').Value = '  +3.69%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '21.05'
$ws.Range('E33').Value = '  +1.98%  '

$ws.Range('E34').Value = '  -0.80%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '156.09'
$ws.Range('E35').Value = '  +0.59%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '4.64'
$ws.Range('E36').Value = '  +2.13%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '6.16'
$ws.Range('E37').Value = '  +2.26%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '27.22'
$ws.Range('E38').Value = '  -1.47%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.26'
$ws.Range('E39').Value = '  -1.26%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0664'
$ws.Range('E40').Value = '  -2.16%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '3.98'
$ws.Range('E41').Value = '  +1.69%  '

$ws.Range('D42').Value = '3.143.76'
$ws.Range('E42').Value = '  +1.10%  '

$ws.Range('B43').Value = 'Stacks'
$ws.Range('C43').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.52'
$ws.Range('E43').Value = '  +9.59%  '

$ws.Range('B44').Value = 'Mantle'
$ws.Range('C44').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.684'
$ws.Range('E44').Value = '  +4.05%  '

$ws.Range('E45').Value = '  +0.01%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.00'
$ws.Range('E46').Value = '  -0.03%  '

$ws.Range('D47').Value = '2.299.74'
$ws.Range('E47').Value = '  -0.18%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0261'
$ws.Range('E48').Value = '  +1.68%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.985'
$ws.Range('E49').Value = '  +4.55%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '20.60'
$ws.Range('E50').Value = '  -1.37%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '6.04'
$ws.Range('E51').Value = '  +1.76%  '
